# Employee training / exam dashboard refresh:
#  - "Period to expire" (H) and "Last update" (I) recalculated against the
#    new reference date (16-Sep-2025 instead of 08-Sep-2025)
#  - Row 19 (LOTO) flips from VALID to NOT VALID and picks up the red
#    "not valid" row styling
#  - Exam Dashboard comments collapse to "date is valid" and the COMMENTS
#    column is narrowed
#  - Header row / title font recolored to white-on-blue

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)   # Exam Dashboard

# ---------------------------------------------------------------------
# 1) Training Dashboard: update "PERIOD TO EXPIRE" (H) + "LAST UPDATE" (I)
#    for every data row (3-20, 22-33). Row 21 is handled separately below
#    because it also changes status/styling.
# ---------------------------------------------------------------------
$rows    = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,24,25,26,27,28,29,30,31,32,33)
$newH    = @(386,382,360,358,405,324,406,342,344,719,446,409,408,327,386,426,427,502,-104,-190,228,-45,182,200,182,199,228,228,354,354)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]

    $hCell = $ws1.Cells.Item($r, 8)
    $hCell.Value = $newH[$i]

    $iCell = $ws1.Cells.Item($r, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value = "16-Sep-2025"
}

# The NumberFormat = "@" trick above forces Excel to keep the date-looking
# text as a literal string instead of silently converting it to a date
# serial, but it also parks each cell on its own one-off style. Re-stamp
# the original (unaffected) formatting on top so every LAST UPDATE cell
# keeps sharing its row's normal style (white rows use style of J3; the
# red "NOT VALID" rows 22/23/25 must keep their own red styling).
$ws1.Range("J3").Copy()
$ws1.Range("I3:I20").PasteSpecial(-4122)
$ws1.Range("I24").PasteSpecial(-4122)
$ws1.Range("I26:I33").PasteSpecial(-4122)

$ws1.Range("J22").Copy()
$ws1.Range("I22").PasteSpecial(-4122)
$ws1.Range("J23").Copy()
$ws1.Range("I23").PasteSpecial(-4122)
$ws1.Range("J25").Copy()
$ws1.Range("I25").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Row 21 (LOTO / SOPs): expired -> pick up the red "NOT VALID" look
#    used elsewhere on the sheet (same formatting as row 22), then update
#    its own values.
# ---------------------------------------------------------------------
$ws1.Range("A22:K22").Copy()
$ws1.Range("A21:K21").PasteSpecial(-4122)

$ws1.Cells.Item(21, 8).Value = 15

$i21 = $ws1.Cells.Item(21, 9)
$i21.NumberFormat = "@"
$i21.Value = "16-Sep-2025"

$ws1.Cells.Item(21, 10).Value = "NOT VALID"

$ws1.Range("J22").Copy()
$ws1.Range("I21").PasteSpecial(-4122)
$ws1.Range("J22").Copy()
$ws1.Range("J21").PasteSpecial(-4122)
$ws1.Cells.Item(21, 10).Value = "NOT VALID"

# ---------------------------------------------------------------------
# 3) Header row + title: white bold text on the dark blue fill.
# ---------------------------------------------------------------------
$ws1.Range("A1").Font.Color = 16777215
$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A1").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

# ---------------------------------------------------------------------
# 4) Exam Dashboard: comments all collapse to "date is valid", and the
#    COMMENTS column narrows from 44 to 15 characters wide.
# ---------------------------------------------------------------------
for ($r = 3; $r -le 9; $r++) {
    $ws2.Cells.Item($r, 5).Value = "date is valid"
}
$ws2.Columns.Item(5).ColumnWidth = 14.17
